$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.446.27"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.081.31"
$ws.Range("E3").Value = "  +1.43%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.05"
$ws.Range("E5").Value = "  +0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.30"
$ws.Range("E6").Value = "  +0.51%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("E8").Value = "  -0.53%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.27"
$ws.Range("E9").Value = "  +1.09%  "

# Row 10
$ws.Range("E10").Value = "  -0.66%  "

# Row 11
$ws.Range("E11").Value = "  -0.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.607.16"
$ws.Range("E12").Value = "  +1.17%  "

# Row 13
$ws.Range("E13").Value = "  +2.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.58"
$ws.Range("E14").Value = "  -5.14%  "

# Row 15
$ws.Range("E15").Value = "  -1.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.521.83"
$ws.Range("E16").Value = "  +0.83%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.077.19"
$ws.Range("E17").Value = "  +1.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.13"
$ws.Range("E18").Value = "  -0.63%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.16"
$ws.Range("E19").Value = "  -1.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.14"
$ws.Range("E20").Value = "  +0.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.88"
$ws.Range("E21").Value = "  +1.06%  "

# Row 22
$ws.Range("E22").Value = "  +0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.501"
$ws.Range("E23").Value = "  -1.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.95"
$ws.Range("E24").Value = "  +1.06%  "

# Row 25
$ws.Range("E25").Value = "  +3.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("E27").Value = "  +2.63%  "

# Row 28
$ws.Range("E28").Value = "  -4.43%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.11"
$ws.Range("E29").Value = "  -1.24%  "

# Row 30
$ws.Range("E30").Value = "  +0.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.83"
$ws.Range("E31").Value = "  +0.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.17"
$ws.Range("E32").Value = "  -3.30%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.79"
$ws.Range("E33").Value = "  +1.54%  "

# Row 34
$ws.Range("B34").Value = "EnergySwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.98"
$ws.Range("E34").Value = "  +11.10%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.54"
$ws.Range("E35").Value = "  -3.08%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.88"
$ws.Range("E36").Value = "  -0.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0675"
$ws.Range("E38").Value = "  -0.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.119.24"
$ws.Range("E39").Value = "  +1.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.76"
$ws.Range("E40").Value = "  -0.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -0.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.671"
$ws.Range("E42").Value = "  +1.66%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.285.69"
$ws.Range("E44").Value = "  +3.91%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0255"
$ws.Range("E45").Value = "  +5.72%  "

# Row 46
$ws.Range("E46").Value = "  -1.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.937"
$ws.Range("E47").Value = "  -0.78%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.98"
$ws.Range("E48").Value = "  -0.49%  "

# Row 49
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.88"
$ws.Range("E49").Value = "  -3.39%  "

# Row 50
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "253.04"
$ws.Range("E50").Value = "  +6.61%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0877"
$ws.Range("E51").Value = "  +1.38%  "
